$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.003.72'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.55%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.303.19'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.19%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.95'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '105.35'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.75%  '
$ws.Range("E7").Value = '  -1.40%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.607'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.66%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.22'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.80%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0913'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.87%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.26'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.90%  '
$ws.Range("E13").Value = '  +0.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.975'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.46%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.53'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.655.97'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.74%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.301.11'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.20%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.131.62'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.63'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.46%  '
$ws.Range("E20").Value = '  -1.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '74.44'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.27%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.46'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.86%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '259.34'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.67%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.29'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.74%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.29'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -7.61%  '
$ws.Range("E26").Value = '  +0.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.97'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.23%  '
$ws.Range("E28").Value = '  +3.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.74'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.95%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.73'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.25%  '
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '164.14'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.59%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0895'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.47%  '
$ws.Range("E33").Value = '  -5.52%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.84'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.94%  '
$ws.Range("E35").Value = '  -2.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.118'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +11.84%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.52'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.02%  '
$ws.Range("E38").Value = '  -1.54%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.77'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.59%  '
$ws.Range("E40").Value = '  -5.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '71.81'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.57%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '98.23'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +8.71%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.46'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.17%  '
$ws.Range("E44").Value = '  -3.16%  '
$ws.Range("E45").Value = '  +0.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.24'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.53%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '112.37'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.02'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.32'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.24%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '74.31'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.93%  '
$ws.Range("E51").Value = '  -0.46%  '
